# Apply 2024-10-31 data refresh to violent-crime-full-year workbook
# Updates year-to-date (2024) totals in column K (and a couple of 2023 column J corrections)
# across the Citywide Totals, By Neighborhood, and per-neighborhood detail sheets.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @(
    @{ Name = "Citywide Totals"; Cells = @{ "K2" = 6780; "K3" = 6998; "J4" = 1840; "K4" = 1446; "K6" = 7652; "J7" = 29306; "K7" = 23376 } },
    @{ Name = "By Neighborhood"; Cells = @{ "K2" = 206; "K5" = 60; "K6" = 167; "K7" = 710; "K8" = 1538; "K11" = 432; "K14" = 115; "K15" = 247; "K19" = 688; "K20" = 565; "K22" = 75; "K23" = 231; "K27" = 219; "K29" = 1267; "K30" = 88; "K33" = 1009; "K35" = 35; "K37" = 796; "K41" = 161; "K42" = 863; "K43" = 188; "K44" = 196; "K46" = 49; "K47" = 156; "K48" = 299; "K49" = 127; "K51" = 290; "K52" = 618; "K54" = 459; "K55" = 251; "J63" = 119; "K63" = 58; "K67" = 909; "K69" = 51; "K76" = 315; "K78" = 270; "K84" = 187; "K85" = 1081; "K88" = 251; "K89" = 348; "K90" = 221; "K92" = 87; "K94" = 312; "K95" = 388; "K99" = 396; "J101" = 29306; "K101" = 23376 } },
    @{ Name = "Bridgeport"; Cells = @{ "K3" = 27; "K7" = 115 } },
    @{ Name = "Auburn Gresham"; Cells = @{ "K2" = 236; "K7" = 710 } },
    @{ Name = "Belmont Cragin"; Cells = @{ "K6" = 142; "K7" = 432 } },
    @{ Name = "Uptown"; Cells = @{ "K6" = 102; "K7" = 348 } },
    @{ Name = "South Shore"; Cells = @{ "K2" = 355; "K3" = 376; "K4" = 56; "K6" = 263; "K7" = 1081 } },
    @{ Name = "Little Village"; Cells = @{ "K6" = 225; "K7" = 618 } },
    @{ Name = "Norwood Park"; Cells = @{ "K6" = 14; "K7" = 51 } },
    @{ Name = "Austin"; Cells = @{ "K2" = 429; "K3" = 467; "K4" = 89; "K6" = 507; "K7" = 1538 } },
    @{ Name = "Garfield Park"; Cells = @{ "K2" = 255; "K3" = 358; "K6" = 318; "K7" = 1009 } },
    @{ Name = "West Pullman"; Cells = @{ "K2" = 133; "K7" = 388 } },
    @{ Name = "Grand Crossing"; Cells = @{ "K2" = 228; "K3" = 262; "K6" = 238; "K7" = 796 } },
    @{ Name = "Woodlawn"; Cells = @{ "K2" = 105; "K3" = 163; "K7" = 396 } },
    @{ Name = "Fuller Park"; Cells = @{ "K2" = 25; "K6" = 31; "K7" = 88 } },
    @{ Name = "North Lawndale"; Cells = @{ "K2" = 247; "K3" = 331; "K6" = 261; "K7" = 909 } },
    @{ Name = "South Deering"; Cells = @{ "K2" = 62; "K7" = 187 } },
    @{ Name = "Lincoln Park"; Cells = @{ "K2" = 29; "K7" = 127 } },
    @{ Name = "Loop"; Cells = @{ "K2" = 73; "K6" = 248; "K7" = 459 } },
    @{ Name = "Englewood"; Cells = @{ "K2" = 358; "K3" = 451; "K7" = 1267 } },
    @{ Name = "Lake View"; Cells = @{ "K3" = 71; "K7" = 299 } },
    @{ Name = "Chatham"; Cells = @{ "K2" = 202; "K3" = 207; "K4" = 33; "K7" = 688 } },
    @{ Name = "Irving Park"; Cells = @{ "K3" = 53; "K7" = 196 } },
    @{ Name = "River North"; Cells = @{ "K3" = 60; "K6" = 159; "K7" = 315 } },
    @{ Name = "Ashburn"; Cells = @{ "K2" = 64; "K7" = 167 } },
    @{ Name = "Hermosa"; Cells = @{ "K3" = 35; "K7" = 161 } },
    @{ Name = "Humboldt Park"; Cells = @{ "K2" = 233; "K3" = 262; "K7" = 863 } },
    @{ Name = "Rogers Park"; Cells = @{ "K2" = 80; "K3" = 69; "K7" = 270 } },
    @{ Name = "Lower West Side"; Cells = @{ "K3" = 74; "K7" = 251 } },
    @{ Name = "Jefferson Park"; Cells = @{ "K3" = 13; "K7" = 49 } },
    @{ Name = "Douglas"; Cells = @{ "K3" = 80; "K4" = 16; "K7" = 231 } },
    @{ Name = "Chicago Lawn"; Cells = @{ "K3" = 182; "K6" = 153; "K7" = 565 } },
    @{ Name = "West Loop"; Cells = @{ "K4" = 25; "K7" = 312 } },
    @{ Name = "Kenwood"; Cells = @{ "K2" = 46; "K3" = 46; "K7" = 156 } },
    @{ Name = "Brighton Park"; Cells = @{ "K2" = 92; "K6" = 73; "K7" = 247 } },
    @{ Name = "Gold Coast"; Cells = @{ "K6" = 21; "K7" = 35 } },
    @{ Name = "Albany Park"; Cells = @{ "K2" = 60; "K7" = 206 } },
    @{ Name = "West Elsdon"; Cells = @{ "K6" = 42; "K7" = 87 } },
    @{ Name = "United Center"; Cells = @{ "K6" = 101; "K7" = 251 } },
    @{ Name = "Armour Square"; Cells = @{ "K3" = 17; "K6" = 27; "K7" = 60 } },
    @{ Name = "Edgewater"; Cells = @{ "K6" = 80; "K7" = 219 } },
    @{ Name = "Washington Heights"; Cells = @{ "K3" = 62; "K6" = 56; "K7" = 221 } },
    @{ Name = "Little Italy, UIC"; Cells = @{ "K6" = 97; "K7" = 290 } },
    @{ Name = "Hyde Park"; Cells = @{ "K3" = 53; "K7" = 188 } },
    @{ Name = "Clearing"; Cells = @{ "K2" = 36; "K3" = 21; "K7" = 75 } }
)

foreach ($update in $sheetUpdates) {
    $ws = $wb.Worksheets($update.Name)
    foreach ($cellRef in $update.Cells.Keys) {
        $ws.Range($cellRef).Value = $update.Cells[$cellRef]
    }
}
